$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 1 (sheet2.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("A4").Value  = "Age at the accident, years"
$ws1.Range("A5").Value  = "Age at the accident, class, years"
$ws1.Range("A8").Value  = "Highest education grade"
$ws1.Range("A9").Value  = "Employment at the accident"
$ws1.Range("A10").Value = "Mountain sport profession"
$ws1.Range("A11").Value = "Search and rescue profession"
$ws1.Range("A17").Value = "Pre-existing diagnosed mental disorder"
$ws1.Range("A18").Value = "Type of pre-existing diagnosed mental disorder"

$b18 = "affective disorder: 2.3% (n = 7)`npersonality disorder: 0.33% (n = 1)`npost-traumatic stress disorder: 0.65% (n = 2)`nsomatoform disorder: 1.6% (n = 5)`nanxiety disorder: 0.65% (n = 2)`nattention-deficit hyperactivity disorder: 0.33% (n = 1)`naddiction: 0.33% (n = 1)`nbulimia nervosa: 0.33% (n = 1)"
$ws1.Range("B18").Value = $b18

# Rows 19-25 (Personality disorder ... Bulimia nervosa) are now folded into B18,
# so remove them entirely; this also shrinks the sheet dimension to A1:B18.
$ws1.Rows("19:25").Delete()

# ---------------------------------------------------------------------------
# Table 2 (sheet3.xml)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")

$ws2.Range("A2").Value  = "Prior mountain sport accidents"
$ws2.Range("A3").Value  = "Mountain sport type"
$ws2.Range("A6").Value  = "Number of injured persons"
$ws2.Range("A7").Value  = "Rescue mode"

$b7 = "self: 50% (n = 155)`ncompanion: 21% (n = 63)`nrescue team: 29% (n = 89)`nn = 307"
$ws2.Range("B7").Value = $b7

$ws2.Range("A10").Value = "Surgical therapy"
$ws2.Range("A11").Value = "Psychological/psychiatric support post accident"
$ws2.Range("A12").Value = "Psychological/psychiatric support need post accident"
$ws2.Range("A13").Value = "Physical health consequences of the accident"
$ws2.Range("A14").Value = "Returned to same mountain sport post accident"
$ws2.Range("A15").Value = "Caution during mountain sport post accident"
$ws2.Range("A16").Value = "Flashback frequency during mountain sport"

# ---------------------------------------------------------------------------
# Table 3 (sheet4.xml)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

$ws3.Range("A4").Value  = "Clinically relevant anxiety symptoms (GAD-7 ≥10)"
$ws3.Range("A6").Value  = "Clinically relevant depression symptoms (PHQ-9 ≥11)"
$ws3.Range("A8").Value  = "Clinically relevant somatizaton symptoms (PHQ-15 ≥11)"
$ws3.Range("A15").Value = "PTSD symptoms (at least one PCL-5 domain positive)"
